$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5.6
$ws.Range("E2").Value = 9.5
$ws.Range("F2").Value = 66

$ws.Range("C3").Value = 6.2
$ws.Range("E3").Value = 9.9
$ws.Range("F3").Value = 64

$ws.Range("C4").Value = 0.5
$ws.Range("E4").Value = 2.6
$ws.Range("F4").Value = 74

$ws.Range("C5").Value = 5.5
$ws.Range("E5").Value = 7.5
$ws.Range("F5").Value = 69

$ws.Range("C6").Value = 6.2
$ws.Range("E6").Value = 14
$ws.Range("F6").Value = 50

$ws.Range("C7").Value = 0.5
$ws.Range("E7").Value = 7.8
$ws.Range("F7").Value = 68

$ws.Range("C8").Value = 4.9
$ws.Range("E8").Value = 7.000000000000001
$ws.Range("F8").Value = 71

$ws.Range("C9").Value = 5.3
$ws.Range("E9").Value = 12.8
$ws.Range("F9").Value = 53

$ws.Range("C10").Value = 5.9
$ws.Range("E10").Value = 10.3
$ws.Range("F10").Value = 62

$ws.Range("C11").Value = 6.2
$ws.Range("E11").Value = 12.5
$ws.Range("F11").Value = 54

$ws.Range("C12").Value = 5.1
$ws.Range("E12").Value = 11.3
$ws.Range("F12").Value = 59

$ws.Range("C13").Value = 5.6
$ws.Range("E13").Value = 12
$ws.Range("F13").Value = 56

$ws.Range("C14").Value = 6.2
$ws.Range("E14").Value = 11.8
$ws.Range("F14").Value = 57

$ws.Range("C15").Value = 6.1
$ws.Range("E15").Value = 7.5
$ws.Range("F15").Value = 69

$ws.Range("C16").Value = 6.1
$ws.Range("E16").Value = 12.2
$ws.Range("F16").Value = 55

$ws.Range("C17").Value = 5.4
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 52

$ws.Range("C18").Value = 6
$ws.Range("E18").Value = 14.2
$ws.Range("F18").Value = 49

$ws.Range("C19").Value = 5.4
$ws.Range("E19").Value = 14.4
$ws.Range("F19").Value = 48

$ws.Range("C20").Value = 5.6
$ws.Range("E20").Value = 13.6
$ws.Range("F20").Value = 51

$ws.Range("C21").Value = 3.9
$ws.Range("E21").Value = 6.3
$ws.Range("F21").Value = 72

$ws.Range("C22").Value = 2.7
$ws.Range("E22").Value = 9.699999999999999
$ws.Range("F22").Value = 65

$ws.Range("C23").Value = 5.7
$ws.Range("E23").Value = 10.3
$ws.Range("F23").Value = 62

$ws.Range("C24").Value = 5.3
$ws.Range("E24").Value = 11.5
$ws.Range("F24").Value = 58

$ws.Range("C25").Value = 6.1
$ws.Range("E25").Value = 15.3
$ws.Range("F25").Value = 46

$ws.Range("C26").Value = 3.6
$ws.Range("E26").Value = 8.200000000000001
$ws.Range("F26").Value = 67

$ws.Range("C27").Value = 5.7
$ws.Range("E27").Value = 16.6
$ws.Range("F27").Value = 44

$ws.Range("C28").Value = 0
$ws.Range("E28").Value = 4.3
$ws.Range("F28").Value = 73

$ws.Range("C29").Value = 6.1
$ws.Range("E29").Value = 15
$ws.Range("F29").Value = 47

$ws.Range("C30").Value = 6
$ws.Range("E30").Value = 17.7
$ws.Range("F30").Value = 38

$ws.Range("C31").Value = 6.2
$ws.Range("E31").Value = 17.6
$ws.Range("F31").Value = 39

$ws.Range("C32").Value = 5.9
$ws.Range("E32").Value = 11.2
$ws.Range("F32").Value = 60

$ws.Range("F33").Value = 32

$ws.Range("C34").Value = 6
$ws.Range("E34").Value = 17.6
$ws.Range("F34").Value = 39

$ws.Range("C35").Value = 6.2
$ws.Range("E35").Value = 19.4
$ws.Range("F35").Value = 26

$ws.Range("C36").Value = 5.6
$ws.Range("E36").Value = 16.3
$ws.Range("F36").Value = 45

$ws.Range("C37").Value = 6.1
$ws.Range("E37").Value = 19.3
$ws.Range("F37").Value = 29

$ws.Range("C38").Value = 6.1
$ws.Range("E38").Value = 18.2
$ws.Range("F38").Value = 37

$ws.Range("C39").Value = 6.2
$ws.Range("E39").Value = 18.5
$ws.Range("F39").Value = 34

$ws.Range("C40").Value = 6.2
$ws.Range("E40").Value = 18.3
$ws.Range("F40").Value = 35

$ws.Range("C41").Value = 6.2
$ws.Range("E41").Value = 19.6
$ws.Range("F41").Value = 25

$ws.Range("C42").Value = 6.2
$ws.Range("E42").Value = 18.3
$ws.Range("F42").Value = 35

$ws.Range("C43").Value = 6
$ws.Range("E43").Value = 19.3
$ws.Range("F43").Value = 27

$ws.Range("C44").Value = 5.3
$ws.Range("E44").Value = 17.5
$ws.Range("F44").Value = 41

$ws.Range("C45").Value = 6.3
$ws.Range("E45").Value = 20.1
$ws.Range("F45").Value = 23

$ws.Range("C46").Value = 6.3
$ws.Range("E46").Value = 20.4
$ws.Range("F46").Value = 22

$ws.Range("F47").Value = 33

$ws.Range("C48").Value = 6.1
$ws.Range("E48").Value = 19.1

$ws.Range("C49").Value = 6.3
$ws.Range("E49").Value = 20.7
$ws.Range("F49").Value = 20

$ws.Range("C50").Value = 5.7
$ws.Range("E50").Value = 17.4

$ws.Range("C51").Value = 6.2
$ws.Range("E51").Value = 20.6
$ws.Range("F51").Value = 21

$ws.Range("C52").Value = 4.9
$ws.Range("E52").Value = 19.2
$ws.Range("F52").Value = 30

$ws.Range("C53").Value = 5.5
$ws.Range("E53").Value = 19.3
$ws.Range("F53").Value = 27

$ws.Range("C54").Value = 6
$ws.Range("E54").Value = 17.2
$ws.Range("F54").Value = 43

$ws.Range("C55").Value = 6.3
$ws.Range("E55").Value = 21.4
$ws.Range("F55").Value = 12

$ws.Range("C56").Value = 6.1
$ws.Range("E56").Value = 21
$ws.Range("F56").Value = 17

$ws.Range("C57").Value = 6.1
$ws.Range("E57").Value = 21.1
$ws.Range("F57").Value = 16

$ws.Range("C58").Value = 6
$ws.Range("E58").Value = 21.3
$ws.Range("F58").Value = 14

$ws.Range("C59").Value = 0
$ws.Range("E59").Value = 10.9
$ws.Range("F59").Value = 61

$ws.Range("C60").Value = 6.2
$ws.Range("E60").Value = 20.9
$ws.Range("F60").Value = 19

$ws.Range("C61").Value = 6.3
$ws.Range("E61").Value = 21.7
$ws.Range("F61").Value = 7

$ws.Range("C62").Value = 6
$ws.Range("E62").Value = 21.6
$ws.Range("F62").Value = 9

$ws.Range("C63").Value = 6.3
$ws.Range("E63").Value = 21.9
$ws.Range("F63").Value = 6

$ws.Range("C64").Value = 6.2
$ws.Range("E64").Value = 22.1
$ws.Range("F64").Value = 5

$ws.Range("C65").Value = 6.1
$ws.Range("E65").Value = 21.7

$ws.Range("C66").Value = 6.2
$ws.Range("E66").Value = 21
$ws.Range("F66").Value = 17

$ws.Range("C67").Value = 6.2
$ws.Range("E67").Value = 21.4
$ws.Range("F67").Value = 12

$ws.Range("C68").Value = 5
$ws.Range("E68").Value = 19.7
$ws.Range("F68").Value = 24

$ws.Range("C69").Value = 6.2
$ws.Range("E69").Value = 22.2

$ws.Range("C70").Value = 6.2
$ws.Range("E70").Value = 22.2
$ws.Range("F70").Value = 3

$ws.Range("C71").Value = 6.3
$ws.Range("E71").Value = 22.3
$ws.Range("F71").Value = 2

$ws.Range("C72").Value = 5.8
$ws.Range("E72").Value = 21.5
$ws.Range("F72").Value = 11

$ws.Range("C73").Value = 6.2
$ws.Range("E73").Value = 21.1
$ws.Range("F73").Value = 15

$ws.Range("C74").Value = 6.3
$ws.Range("E74").Value = 22.9
$ws.Range("F74").Value = 1

$ws.Range("C75").Value = 5.8
$ws.Range("E75").Value = 21.6
$ws.Range("F75").Value = 9
